$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.066.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.17%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.635.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.15"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5233"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06295"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.09%  "
$ws.Range("E10").Value = "  -2.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07653"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.41%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.627.74"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.55%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.414"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.859.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.93%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5503"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8166"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.08%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.90"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.056.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.682"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "188.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.08%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.15"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.79%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.140"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.43%  "
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.62"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.83%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1211"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.85%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.389"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.76"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.399"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05950"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -5.02%  "
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.441"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.05%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.404"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.26%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.634"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9831"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.396"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.760"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.93%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -5.39%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01615"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8522"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.14%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.733"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -6.28%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.033.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.785.95"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.83%  "
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "55.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.76%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0₈105"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.26%  "
$ws.Range("E48").Value = "  -0.26%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.024"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05165"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4219"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.62%  "

Write-Host "Applied all changes"
